$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting the existing rows 13-17 down to 14-18
$ws.Range("A13:R13").Insert()

# Populate the newly inserted row 13 with its data
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C13").Value = "Arica y Parinacota"
$ws.Range("D13").Value = 44658
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E13").Value = 15
$ws.Range("F13").Value = 100112006
$ws.Range("G13").Value = "Repollo"
$ws.Range("H13").Value = "Copenhague"
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 600
$ws.Range("L13").Value = 650
$ws.Range("M13").Value = 625
$ws.Range("N13").Value = "$/unidad"
$ws.Range("O13").Value = "Región de Arica y Parinacota"
$ws.Range("P13").Value = 625
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"
